$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.826.34"
$ws.Range("E2").Value = "  +3.64%  "
$ws.Range("D3").Value = "3.684.91"
$ws.Range("E3").Value = "  +8.77%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "588.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.05%  "
$ws.Range("D7").Value = "3.677.41"
$ws.Range("E7").Value = "  +8.78%  "
$ws.Range("E8").Value = "  +5.10%  "
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.202"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.97%  "
$ws.Range("E11").Value = "  +4.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.96"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.33%  "
$ws.Range("E13").Value = "  +1.45%  "
$ws.Range("D14").Value = "4.281.30"
$ws.Range("E14").Value = "  +8.88%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "683.48"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.12%  "
$ws.Range("E16").Value = "  +4.41%  "
$ws.Range("D17").Value = "71.894.68"
$ws.Range("E17").Value = "  +3.51%  "
$ws.Range("D18").Value = "3.679.40"
$ws.Range("E18").Value = "  +8.57%  "
$ws.Range("E19").Value = "  +2.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.15%  "
$ws.Range("E22").Value = "  +3.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.17"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +15.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "17.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.03%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "103.51"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.35%  "
$ws.Range("E26").Value = "  +3.87%  "
$ws.Range("E27").Value = "  +5.58%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.18"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "35.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +6.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.34%  "
$ws.Range("E31").Value = "  +6.30%  "
$ws.Range("E32").Value = "  +9.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "575.44"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.66%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.31"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.52%  "
$ws.Range("E35").Value = "  +4.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.68"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.19%  "
$ws.Range("D37").Value = "3.747.99"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.999"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("E39").Value = "  +3.27%  "
$ws.Range("E40").Value = "  +4.24%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "35.44"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("E42").Value = "  +4.93%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0462"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.82%  "
$ws.Range("E44").Value = "  +2.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.349"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.02%  "
$ws.Range("E46").Value = "  +7.92%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E48").Value = "  +4.00%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.999"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "134.15"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.52%  "
